# #5: cash & deposit done
# Rebuild the "存款" (deposit) sheet so every row carries the full
# property-record column set (bank / deposit_type / currency / owner /
# total / property_category / category / date / legislator_name /
# legislator_id / source_file / index), matching the layout already used
# on the other property sheets (土地/建物/汽車).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- header row -----------------------------------------------------
# columns H:M are brand new - carry over the bold/bordered header look
# already used for B1:G1 so the new cells format the same way.
$ws.Range("G1").Copy()
$ws.Range("H1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# ---- data rows --------------------------------------------------------
# index, bank, deposit_type, currency, owner, total
$rows = @(
    @(44, "臺灣銀行城中分行",         "活期存款", "美金",   "邱議瑩", 59440),
    @(45, "臺灣銀行群賢分行",         "活期存款", "新臺幣", "邱議瑩", 5905403),
    @(46, "臺灣新光商業銀行屏東分行", "定期存款", "新臺幣", "邱議瑩", 300000),
    @(47, "臺灣新光商業銀行屏東分行", "活期存款", "新臺幣", "邱議瑩", 1800000),
    @(48, "彰化商業銀行屏東分行",     "活期存款", "新臺幣", "邱議瑩", 2835),
    @(49, "合作金庫商業銀行營業部",   "活期存款", "新臺幣", "邱議瑩", 2057),
    @(50, "元大商業銀行營業部",       "活期存款", "新臺幣", "邱議瑩", 102866),
    @(51, "京城商業銀行營業部",       "活期存款", "新臺幣", "邱議瑩", 53059)
)

$r = 2
foreach ($row in $rows) {
    $idx = $row[0]

    # new H:M cells on this row inherit the plain data-row look already
    # used on A:G (thin border / normal weight)
    $ws.Range("G" + $r).Copy()
    $ws.Range("H" + $r + ":M" + $r).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    # force the register date to stay plain text (matches the other
    # property sheets) instead of being auto-parsed into a date serial
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 9).Value = "2011-11-25"
    $ws.Cells.Item($r, 10).Value = "邱議瑩"
    $ws.Cells.Item($r, 11).Value = 913
    $ws.Cells.Item($r, 12).Value = "tmpab161"
    $ws.Cells.Item($r, 13).Value = $idx
    $r = $r + 1
}
